$wb = $excel.ActiveWorkbook

# Grab references to the sheets we need, up front (before any renames happen)
$wsWawancara    = $wb.Worksheets.Item("Wawancara")
$wsSeleksi      = $wb.Worksheets.Item("Seleksi")
$wsPenawaranOld = $wb.Worksheets.Item("Penawaran")
$wsPelamar      = $wb.Worksheets.Item("Pelamar")
$wsLowongan     = $wb.Worksheets.Item("Lowongan")

# 1. Delete the old "Penawaran" sheet entirely - it is being replaced
[void]$wsPenawaranOld.Delete()

# 2. The old "Wawancara" sheet becomes the new "Penawaran" sheet (reusing its slot)
$wsWawancara.Name = "Penawaran"

# 3. The old "Seleksi" sheet becomes the new "Wawancara" sheet (reusing its slot)
$wsSeleksi.Name = "Wawancara"

# --- Rebuild the (new) "Penawaran" sheet content ---
$penawaran = $wsWawancara
$penawaran.Rows("1:10").Delete()
$penawaran.Columns("C:Z").Delete()

$penawaran.Range("A1").Value = "Kode Pekerjaan"
$penawaran.Range("B1").Value = "Benefit"
$penawaran.Range("A2").Value = "L001"
$penawaran.Range("B2").Value = "Asuransi Kesehatan"
$penawaran.Range("A3").Value = "L001"
$penawaran.Range("B3").Value = "Uang Transport"

$penawaran.Rows("1:3").RowHeight = 18.75

# --- Rebuild the (new) "Wawancara" sheet content ---
# (the old "Seleksi" data only used columns A:D, so no column trimming is
# needed here - we're instead extending it out to F)
$wawancara = $wsSeleksi
$wawancara.Rows("1:10").Delete()

$wawancara.Range("A1").Value = "Kode Lowongan"
$wawancara.Range("B1").Value = "Kode Pelamar"
$wawancara.Range("C1").Value = "Nama Pelamar"
$wawancara.Range("D1").Value = "Posisi"
$wawancara.Range("E1").Value = "Jadwal Tanggal"
$wawancara.Range("F1").Value = "Jadwal Jam"

$wawancara.Range("A2").Value = "L003"
$wawancara.Range("B2").Value = "P003"
$wawancara.Range("C2").Value = "Jidan"
$wawancara.Range("D2").Value = "Proggrammer"

# Force these two as plain text so they aren't auto-converted into a date
# serial number / numeric value by Excel's type inference.
$wawancara.Range("E2:F2").NumberFormat = "@"
$wawancara.Range("E2").Value = "2024-07-03"
$wawancara.Range("F2").Value = "18.00"

# --- Update "Pelamar" sheet: interview status column ---
$wsPelamar.Range("E2").Value = "Belum"
$wsPelamar.Range("E3").Value = "Belum"
$wsPelamar.Range("E4").Value = "Proses"

# Pelamar becomes the active tab with E2 selected
[void]$wsPelamar.Activate()
[void]$wsPelamar.Range("E2").Select()

# --- Update "Lowongan" sheet: row 5 gets the standard row height ---
$wsLowongan.Rows(5).RowHeight = 18.75
